$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "Campaign Dates that use Taurus constellation 2022: January 16-25",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " 2022 Campaign Dates that use Taurus constellation: January 16-25",
    2
)
